# Update the "dSF" column (F) values per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    3  = -2
    6  = -4
    7  = 4
    8  = 2
    10 = -4
    11 = -4
    12 = 2
    13 = -6
    14 = -2
    15 = -2
    16 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
